$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66 ("2025-05-05") previously carried "NA" in column C. The script
# run that produced this update determined there *is* a term to report for
# 2025-05-06, so the "NA" placeholder moves down to a brand-new row for
# that date, and row 66's C cell goes back to blank (consistent with every
# other "Rien ne nous concerne aujourd'hui !" row).
$ws.Range("C66").ClearContents()

# Append the new row 67 with the fresh date's data.
# The leading apostrophe keeps the date a literal text value (matching the
# rest of column A) instead of Excel auto-converting it to a date serial.
$ws.Range("A67").Value = "'2025-05-06"
$ws.Range("A67").Style = "Normal"
$ws.Range("B67").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C67").Value = "NA"
$ws.Range("D67").Value = 1
